# Auto-generated Excel COM-interop script to apply market-data refresh diff
# to Seraph_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 284.16666
$ws.Range("I33").Value = 176.25
$ws.Range("K33").Value = 176.25
$ws.Range("M33").Value = 52.75
$ws.Range("H80").Value = 544.5714
$ws.Range("I80").Value = 408.3846
$ws.Range("J80").Value = 765.875
$ws.Range("K80").Value = 1225.1538
$ws.Range("L80").Value = 2297.625
$ws.Range("M80").Value = -227.1538
$ws.Range("N80").Value = -4293.625
$ws.Range("H83").Value = 544.5714
$ws.Range("I83").Value = 408.3846
$ws.Range("J83").Value = 765.875
$ws.Range("K83").Value = 3675.4614
$ws.Range("L83").Value = 6892.875
$ws.Range("M83").Value = 1316.5386
$ws.Range("N83").Value = -16876.875
$ws.Range("H100").Value = 1789
$ws.Range("I100").Value = 1535.375
$ws.Range("J100").Value = 2465.3333
$ws.Range("K100").Value = 1535.375
$ws.Range("L100").Value = 2465.3333
$ws.Range("M100").Value = -994.375
$ws.Range("N100").Value = -3547.3333
$ws.Range("H101").Value = 470.85715
$ws.Range("I101").Value = 470.85715
$ws.Range("K101").Value = 1412.57145
$ws.Range("M101").Value = 209.4285500000001
$ws.Range("H125").Value = 969.6
$ws.Range("I125").Value = 968.44446
$ws.Range("K125").Value = 8716.00014
$ws.Range("M125").Value = -6256.00014
$ws.Range("H132").Value = 1967.75
$ws.Range("I132").Value = 1955.7273
$ws.Range("J132").Value = 2100
$ws.Range("K132").Value = 5867.1819
$ws.Range("L132").Value = 6300
$ws.Range("M132").Value = -3337.1819
$ws.Range("N132").Value = -11360

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2524.1538
$ws.Range("I45").Value = 2526.1667
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 2526.1667
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -2149.1667
$ws.Range("N45").Value = -3254
$ws.Range("H74").Value = 925
$ws.Range("I74").Value = 893.05554
$ws.Range("K74").Value = 893.05554
$ws.Range("M74").Value = -19.05553999999995
$ws.Range("H77").Value = 925
$ws.Range("I77").Value = 893.05554
$ws.Range("K77").Value = 4465.2777
$ws.Range("M77").Value = -97.27769999999964
$ws.Range("H132").Value = 1133.9375
$ws.Range("I132").Value = 1050.3846
$ws.Range("K132").Value = 3151.1538
$ws.Range("M132").Value = -621.1538

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 35944.5
$ws.Range("I96").Value = 7890
$ws.Range("J96").Value = 63999
$ws.Range("K96").Value = 7890
$ws.Range("L96").Value = 63999
$ws.Range("M96").Value = -5144
$ws.Range("N96").Value = -69491
$ws.Range("H99").Value = 2856.3635
$ws.Range("I99").Value = 3010
$ws.Range("J99").Value = 2822.2222
$ws.Range("K99").Value = 3010
$ws.Range("L99").Value = 2822.2222
$ws.Range("M99").Value = -1512
$ws.Range("N99").Value = -5818.2222
$ws.Range("H107").Value = 3136.5
$ws.Range("I107").Value = 3087.5
$ws.Range("K107").Value = 3087.5
$ws.Range("M107").Value = -1167.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 275
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150
$ws.Range("H62").Value = 60443.43
$ws.Range("I62").Value = 3035
$ws.Range("K62").Value = 3035
$ws.Range("M62").Value = -2411
$ws.Range("H65").Value = 60443.43
$ws.Range("I65").Value = 3035
$ws.Range("K65").Value = 15175
$ws.Range("M65").Value = -12055
$ws.Range("H80").Value = 39999.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 39999.5
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 39999.5
$ws.Range("N80").Value = -42245.5
$ws.Range("H83").Value = 39999.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 39999.5
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 119998.5
$ws.Range("N83").Value = -131230.5
$ws.Range("H141").Value = 139552.8
$ws.Range("J141").Value = 139552.8
$ws.Range("L141").Value = 139552.8
$ws.Range("N141").Value = -149912.8

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 142.33333
$ws.Range("J23").Value = 323
$ws.Range("L23").Value = 969
$ws.Range("N23").Value = -1439
$ws.Range("H38").Value = 7824.923
$ws.Range("I38").Value = 11278.223
$ws.Range("K38").Value = 33834.669
$ws.Range("M38").Value = -33487.669

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 45000
$ws.Range("J15").Value = 45000
$ws.Range("L15").Value = 45000
$ws.Range("N15").Value = -45576
$ws.Range("H80").Value = 7033.1113
$ws.Range("I80").Value = 4950
$ws.Range("K80").Value = 4950
$ws.Range("M80").Value = -3952
$ws.Range("H81").Value = 45000
$ws.Range("J81").Value = 45000
$ws.Range("L81").Value = 45000
$ws.Range("N81").Value = -46996
$ws.Range("H83").Value = 7033.1113
$ws.Range("I83").Value = 4950
$ws.Range("K83").Value = 24750
$ws.Range("M83").Value = -19758
$ws.Range("H84").Value = 45000
$ws.Range("J84").Value = 45000
$ws.Range("L84").Value = 135000
$ws.Range("N84").Value = -144984
$ws.Range("H132").Value = 1376.3077
$ws.Range("I132").Value = 814.2
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 2442.6
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = 87.39999999999964
$ws.Range("N132").Value = -14810

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2608
$ws.Range("I40").Value = 2608
$ws.Range("K40").Value = 2608
$ws.Range("M40").Value = -2472
$ws.Range("H43").Value = 3390999.2
$ws.Range("I43").Value = 7998.5
$ws.Range("K43").Value = 7998.5
$ws.Range("M43").Value = -7805.5
$ws.Range("H68").Value = 3601.5
$ws.Range("I68").Value = 2200
$ws.Range("K68").Value = 2200
$ws.Range("M68").Value = -1451
$ws.Range("H71").Value = 3601.5
$ws.Range("I71").Value = 2200
$ws.Range("K71").Value = 11000
$ws.Range("M71").Value = -7256
$ws.Range("H82").Value = 2407.7273
$ws.Range("I82").Value = 3072
$ws.Range("K82").Value = 3072
$ws.Range("M82").Value = -2711
$ws.Range("H85").Value = 2407.7273
$ws.Range("I85").Value = 3072
$ws.Range("K85").Value = 3072
$ws.Range("M85").Value = -1824
$ws.Range("H100").Value = 1560.4
$ws.Range("I100").Value = 1560.4
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1560.4
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -1019.4
$ws.Range("H106").Value = 17099.8
$ws.Range("J106").Value = 17099.8
$ws.Range("L106").Value = 17099.8
$ws.Range("N106").Value = -19623.8
$ws.Range("H136").Value = 7058.375
$ws.Range("I136").Value = 6613.5
$ws.Range("K136").Value = 19840.5
$ws.Range("M136").Value = -17290.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 16000
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H68").Value = 52999
$ws.Range("J68").Value = 52999
$ws.Range("L68").Value = 52999
$ws.Range("N68").Value = -54621
$ws.Range("H71").Value = 52999
$ws.Range("J71").Value = 52999
$ws.Range("L71").Value = 158997
$ws.Range("N71").Value = -167109
$ws.Range("H99").Value = 38999
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H132").Value = 2846.853
$ws.Range("I132").Value = 2503.423
$ws.Range("K132").Value = 7510.268999999999
$ws.Range("M132").Value = -4980.268999999999
$ws.Range("H136").Value = 1844.7858
$ws.Range("I136").Value = 1032.7
$ws.Range("K136").Value = 3098.1
$ws.Range("M136").Value = -548.1000000000004
